# Apply the "Add files via upload" edit: insert a new "affiliate marketing"
# row after the existing "blockchain" row, replace what is now the old
# "motivation quotes" row that slid into row 17 with another affiliate
# marketing row, and append two brand-new rows at the bottom of the table.
#
# NOTE: the order in which brand-new text is typed in matters, because it
# controls the order new entries land in the shared-strings table -- so we
# deliberately type the "הכנסה פסיבית" row before any "affiliate marketing"
# text, matching how the saved workbook ended up ordering its new strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new (currently blank) row before row 5; this shifts old rows
# 5-20 down to 6-21 without touching any cell contents yet.
$ws.Rows.Item(5).Insert()

# Append the two brand-new trailing rows first, so "הכנסה פסיבית" is the
# first new shared string introduced by this edit.
$ws.Cells.Item(22, 1).Value = "הכנסה פסיבית"
$ws.Cells.Item(22, 2).Value = "passive.income.nadi.myfirstdrawermenuproject2"

$ws.Cells.Item(23, 1).Value = "affiliate marketing"
$ws.Cells.Item(23, 2).Value = "passive.income.nadi.affiliatemarketingforpassiveincome"

# Populate the newly inserted row 5 with the affiliate-marketing entry.
$ws.Cells.Item(5, 1).Value = "affiliate marketing"
$ws.Cells.Item(5, 2).Value = "passive.income.nadi.affiliatemarketingforpassiveincome"

# After the shift, the row that used to be row 16 ("motivation quotes") is
# now row 17; overwrite it with another affiliate-marketing row.
$ws.Cells.Item(17, 1).Value = "affiliate marketing"
$ws.Cells.Item(17, 2).Value = "passive.income.nadi.affiliatemarketingforpassiveincome"

# Match the refreshed selection state recorded in the saved workbook: the
# selection moves to the newly inserted row.
$ws.Range("A5:B5").Select()
